# Update the cached text of the "datetimeFigureOut" date placeholder fields
# across the Slide Master, every Slide Layout, and the Notes Master:
# "12/2/2020" -> "12/4/2020"

$p = $ppt.ActivePresentation

function Set-DateFieldText {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "12/2/2020") {
                $tr.Text = "12/4/2020"
            }
        }
    }
}

# Slide Master
Set-DateFieldText $p.SlideMaster.Shapes

# Every Slide Layout under the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Set-DateFieldText $layouts.Item($li).Shapes
}

# Notes Master
Set-DateFieldText $p.NotesMaster.Shapes
